# Updated cryptos list values (price + 1h volume %) per upstream diff.
# Column D ("Price") values that look like plain decimal numbers get
# coerced to numeric cells by Excel on assignment, which would lose the
# original text formatting (e.g. "11.00" -> 11, "0.08019" -> 0.08019 as a
# number instead of text). Force those through as text, matching the
# original inlineStr/shared-string cell type, then restore the default
# "Normal" style so no stray number-format styling is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

$ws.Range("D2").Value = '27.995.92'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '1.890.59'
$ws.Range("E3").Value = '  +0.53%  '
$ws.Range("E4").Value = '  +1.40%  '
Set-TextValue 'D5' '336.05'
$ws.Range("E5").Value = '  +0.88%  '
$ws.Range("E6").Value = '  +1.25%  '
Set-TextValue 'D7' '0.4709'
$ws.Range("E7").Value = '  -0.59%  '
Set-TextValue 'D8' '0.3952'
$ws.Range("E8").Value = '  -0.53%  '
Set-TextValue 'D9' '46.83'
$ws.Range("E9").Value = '  -2.88%  '
Set-TextValue 'D10' '0.08019'
$ws.Range("E10").Value = '  -0.47%  '
Set-TextValue 'D11' '1.019'
$ws.Range("E11").Value = '  -0.90%  '
Set-TextValue 'D12' '21.79'
$ws.Range("E12").Value = '  -0.44%  '
$ws.Range("D13").Value = '1.891.47'
$ws.Range("E13").Value = '  -0.12%  '
Set-TextValue 'D14' '5.996'
$ws.Range("E14").Value = '  +0.52%  '
Set-TextValue 'D15' '7.174'
$ws.Range("E15").Value = '  -0.37%  '
Set-TextValue 'D16' '1.019'
$ws.Range("E16").Value = '  +1.33%  '
Set-TextValue 'D17' '0.06796'
$ws.Range("E17").Value = '  +2.65%  '
Set-TextValue 'D18' '87.98'
$ws.Range("E18").Value = '  +0.77%  '
Set-TextValue 'D19' '0.00001052'
$ws.Range("E19").Value = '  +0.08%  '
Set-TextValue 'D20' '17.21'
$ws.Range("E20").Value = '  -0.56%  '
$ws.Range("E21").Value = '  +1.28%  '
$ws.Range("D22").Value = '27.993.85'
$ws.Range("E22").Value = '  +0.21%  '
Set-TextValue 'D23' '5.506'
$ws.Range("E23").Value = '  -0.26%  '
Set-TextValue 'D24' '11.00'
$ws.Range("E24").Value = '  -0.68%  '
Set-TextValue 'D25' '2.365'
$ws.Range("E25").Value = '  +2.14%  '
$ws.Range("D26").Value = '2.114.27'
$ws.Range("E26").Value = '  -0.27%  '
Set-TextValue 'D27' '159.43'
$ws.Range("E27").Value = '  +1.05%  '
$ws.Range("E28").Value = '  -1.10%  '
Set-TextValue 'D29' '2.103'
$ws.Range("E29").Value = '  -0.32%  '
Set-TextValue 'D30' '5.496'
$ws.Range("E30").Value = '  -2.31%  '
Set-TextValue 'D31' '121.69'
$ws.Range("E31").Value = '  -0.78%  '
Set-TextValue 'D32' '0.09574'
$ws.Range("E32").Value = '  -0.02%  '
Set-TextValue 'D33' '0.9664'
$ws.Range("E33").Value = '  -1.98%  '
Set-TextValue 'D34' '3.648'
$ws.Range("E34").Value = '  +0.67%  '
$ws.Range("E35").Value = '  +0.55%  '
$ws.Range("E36").Value = '  -7.00%  '
Set-TextValue 'D37' '0.06131'
$ws.Range("E37").Value = '  -0.01%  '
Set-TextValue 'D38' '0.02251'
$ws.Range("E38").Value = '  -0.62%  '
Set-TextValue 'D39' '1.217'
$ws.Range("E39").Value = '  -1.03%  '
Set-TextValue 'D40' '8.241'
$ws.Range("E40").Value = '  -0.17%  '
Set-TextValue 'D41' '0.5965'
Set-TextValue 'D42' '0.1905'
$ws.Range("E42").Value = '  -0.36%  '
$ws.Range("E43").Value = '  -0.11%  '
$ws.Range("E44").Value = '  -0.36%  '
Set-TextValue 'D45' '0.5708'
$ws.Range("E45").Value = '  -0.24%  '
Set-TextValue 'D46' '12.28'
$ws.Range("E46").Value = '  +0.03%  '
$ws.Range("B47").Value = 'PancakeSwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-TextValue 'D47' '3.409'
$ws.Range("E47").Value = '  -0.21%  '
$ws.Range("B48").Value = 'NEARProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue 'D48' '1.946'
$ws.Range("E48").Value = '  -0.24%  '
Set-TextValue 'D49' '0.06874'
$ws.Range("E49").Value = '  +0.58%  '
Set-TextValue 'D50' '113.76'
$ws.Range("E50").Value = '  -0.17%  '
Set-TextValue 'D51' '1.071'
$ws.Range("E51").Value = '  -0.53%  '
